# Generate Report for Handoff
#
# The "ff945447-2701-4791-9b2d-41e05a4160a7" entity has now been handed
# off again (newer handoff timestamp, status "Ready for handoff"), while
# "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb" keeps the data it already had
# (it is now reported as still "Handed back: in sync with en-US").
# Concretely the two rows (2 and 3) on every sheet swap identity, and
# the row that ends up in row 3 (now 0a3ca150) picks up a fresh status /
# timestamp.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "ff945447-2701-4791-9b2d-41e05a4160a7.md"
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"
$ov.Range("D2").Value = "2016-03-21 10:42:44"

$ov.Range("A3").Value = "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = "2016-03-21 10:43:47"

# Hyperlinks keep the very same target addresses they already had, only
# the displayed label needs to swap along with the cell text. Rebuild
# the hyperlinks collection (in the same order) so relationship ids stay
# stable (rId2, rId3, ...).
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/29711e7fdf5ec62afc8d72259ecba2f945528085/e2e/0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md", "", "", "ff945447-2701-4791-9b2d-41e05a4160a7.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/29711e7fdf5ec62afc8d72259ecba2f945528085/e2e/ff945447-2701-4791-9b2d-41e05a4160a7.md", "", "", "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md")

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "ff945447-2701-4791-9b2d-41e05a4160a7.md"
$zh.Range("B2").Value = ".md"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("D2").Value = "ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.zh-cn.xlf"
$zh.Range("E2").Value = "2016-03-21 10:42:39"
$zh.Range("F2").Value = "ff945447-2701-4791-9b2d-41e05a4160a7.md"
$zh.Range("G2").Value = "ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.zh-cn.xlf"
$zh.Range("H2").Value = "2016-03-21 10:43:09"
$zh.Range("J2").Value = "Include"

$zh.Range("A3").Value = "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("D3").Value = "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.zh-cn.xlf"
$zh.Range("E3").Value = "2016-03-21 10:43:44"
$zh.Range("F3").Value = "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md"
$zh.Range("G3").Value = "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.zh-cn.xlf"
$zh.Range("H3").Value = "2016-03-21 10:43:09"
$zh.Range("J3").Value = "Include"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/29711e7fdf5ec62afc8d72259ecba2f945528085/e2e/0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md", "", "", "ff945447-2701-4791-9b2d-41e05a4160a7.md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/62c9f873c35295e4459970c75ca8ba3939a48afe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.zh-cn.xlf", "", "", "ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/9f26811d20a6e1d97b7232eecc83a26fb76876dc/e2e/0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md", "", "", "ff945447-2701-4791-9b2d-41e05a4160a7.md")
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bd913ae943a8d0699891951882711a4f58a7b3d8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.zh-cn.xlf", "", "", "ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/29711e7fdf5ec62afc8d72259ecba2f945528085/e2e/ff945447-2701-4791-9b2d-41e05a4160a7.md", "", "", "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md")
$zh.Hyperlinks.Add($zh.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/62c9f873c35295e4459970c75ca8ba3939a48afe/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.zh-cn.xlf", "", "", "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/9f26811d20a6e1d97b7232eecc83a26fb76876dc/e2e/ff945447-2701-4791-9b2d-41e05a4160a7.md", "", "", "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md")
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bd913ae943a8d0699891951882711a4f58a7b3d8/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.zh-cn.xlf", "", "", "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "ff945447-2701-4791-9b2d-41e05a4160a7.md"
$de.Range("B2").Value = ".md"
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("D2").Value = "ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.de-de.xlf"
$de.Range("E2").Value = "2016-03-21 10:42:44"
$de.Range("F2").Value = "ff945447-2701-4791-9b2d-41e05a4160a7.md"
$de.Range("G2").Value = "ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.de-de.xlf"
$de.Range("H2").Value = "2016-03-21 10:43:17"
$de.Range("J2").Value = "Include"

$de.Range("A3").Value = "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("D3").Value = "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.de-de.xlf"
$de.Range("E3").Value = "2016-03-21 10:43:47"
$de.Range("F3").Value = "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md"
$de.Range("G3").Value = "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.de-de.xlf"
$de.Range("H3").Value = "2016-03-21 10:43:17"
$de.Range("J3").Value = "Include"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/29711e7fdf5ec62afc8d72259ecba2f945528085/e2e/0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md", "", "", "ff945447-2701-4791-9b2d-41e05a4160a7.md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8410f360bd6907985443a2555ca5e17d37425994/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.de-de.xlf", "", "", "ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/aafc33ec9209fc0be215a1c3f1a7d60e0e485e5b/e2e/0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md", "", "", "ff945447-2701-4791-9b2d-41e05a4160a7.md")
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/08bdd71a5fbfd31e27f319f7cb2f215ecc8842ee/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.de-de.xlf", "", "", "ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/29711e7fdf5ec62afc8d72259ecba2f945528085/e2e/ff945447-2701-4791-9b2d-41e05a4160a7.md", "", "", "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md")
$de.Hyperlinks.Add($de.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8410f360bd6907985443a2555ca5e17d37425994/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.de-de.xlf", "", "", "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/aafc33ec9209fc0be215a1c3f1a7d60e0e485e5b/e2e/ff945447-2701-4791-9b2d-41e05a4160a7.md", "", "", "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.md")
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/08bdd71a5fbfd31e27f319f7cb2f215ecc8842ee/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/ff945447-2701-4791-9b2d-41e05a4160a7.93c8d1b8c9ab71c0f1001f53f2d8d805e7efbae6.de-de.xlf", "", "", "0a3ca150-d4e6-4814-ab2e-bedcf42e27fb.7d6e6f73e86ad0afe8eb81a66213da2cba5f59c6.de-de.xlf")
